$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# GLOBAL RESULTS sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("GLOBAL RESULTS")

$ws.Range("C6").Value = 22163.090959821988
$ws.Range("C7").Value = 21783.090959821988
$ws.Range("C8").Value = 21498.198231027323
$ws.Range("C12").Value = 3064.9572025323923
$ws.Range("C13").Value = 306.05819399999996
$ws.Range("C14").Value = 19098.133757289597
$ws.Range("C15").Value = 18718.133757289597
$ws.Range("C16").Value = 12258.1337572896
$ws.Range("C17").Value = 11952.0755632896
$ws.Range("C18").Value = 11331.651563289603
$ws.Range("C23").Value = 217345.67596113821
$ws.Range("C24").Value = 213619.1489611382
$ws.Range("C25").Value = 210825.30568230402
$ws.Range("C29").Value = 3001.4055881900995
$ws.Range("C30").Value = 187288.71341092396
$ws.Range("C31").Value = 183562.18641092395
$ws.Range("C32").Value = 120211.22741092404
$ws.Range("C33").Value = 117209.82182273391
$ws.Range("C34").Value = 111125.54080313396

# ---------------------------------------------------------------
# FUSELAGE sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("FUSELAGE")

$ws.Range("C7").Value = 2409.0
$ws.Range("D7").Value = -7.581642929074921
$ws.Range("C8").Value = 2844.0
$ws.Range("D8").Value = 9.106603366422135
$ws.Range("C9").Value = 2588.0
$ws.Range("D9").Value = -0.7145254879393506
$ws.Range("C12").Value = 2940.0
$ws.Range("D12").Value = 12.789526686807678

# ---------------------------------------------------------------
# WING sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("WING")

$ws.Range("C7").Value = 2312.0
$ws.Range("D7").Value = 33.04560494892824
$ws.Range("C13").Value = 1904.7142857142858
$ws.Range("D13").Value = 9.608072836385334

# ---------------------------------------------------------------
# LANDING GEARS sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("LANDING GEARS")

$ws.Range("C5").Value = 740.0
$ws.Range("D5").Value = 6.459502229895028
$ws.Range("C6").Value = 887.0
$ws.Range("D6").Value = 27.607538483671473
$ws.Range("C7").Value = 1007.0
$ws.Range("D7").Value = 44.87124154797878
$ws.Range("C8").Value = 898.0
$ws.Range("D8").Value = 29.190044597899643
$ws.Range("C9").Value = 883.0
$ws.Range("D9").Value = 27.032081714861196
